# "correção da planilha e sprints"
# Adds a new "RESULTADO EM HORAS" column to the Tabela1 table (between
# "QUEM REALIZOU" and "STATUS"), fills it with per-task duration strings,
# removes the trailing blank row, and widens the new column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$lo = $ws.ListObjects.Item(1)

# ---------------------------------------------------------------------
# 1) Drop the last (blank) separator row of the table (old row 40).
# ---------------------------------------------------------------------
$ws.Rows.Item(40).Delete()

# ---------------------------------------------------------------------
# 2) Insert a new column in front of the STATUS column (old F), pushing
#    STATUS + its data one column to the right (F -> G).
# ---------------------------------------------------------------------
$ws.Range("F1").EntireColumn.Insert()

# Resize the table so it covers the new column and the shrunk row range.
$lo.Resize($ws.Range("B2:G39"))

# ---------------------------------------------------------------------
# 3) Headers.
# ---------------------------------------------------------------------
$ws.Range("F2").Value = "RESULTADO EM HORAS"
$ws.Range("G2").Value = "STATUS"

# ---------------------------------------------------------------------
# 4) A handful of rows inherited the neighbouring (non-default) cell
#    style when the column was inserted; the real data rows use the
#    plain style carried by F3 (style id 2), so re-apply that format.
# ---------------------------------------------------------------------
$ws.Range("F3").Copy()
$ws.Range("F25").PasteSpecial(-4122) | Out-Null
$ws.Range("F27").PasteSpecial(-4122) | Out-Null
$ws.Range("F30").PasteSpecial(-4122) | Out-Null
$ws.Range("F37").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------
# 5) Fill in "RESULTADO EM HORAS" values for every task row (separator
#    rows 5, 9, 13, 17, 21, 26, 32 stay blank).
# ---------------------------------------------------------------------
$values = [ordered]@{
    3  = " 2 Horas e 20 Minutos"
    4  = "2 Horas e 15 Minutos"
    6  = "2 Horas 30 Minutos"
    7  = "2 Horas e 15 Minutos"
    8  = "1 Hora e 15 Minutos"
    10 = "1 Hora"
    11 = "1 Hora"
    12 = "1 Hora"
    14 = "1 Hora e 15 Minutos"
    15 = "1 Hora e 25 Minutos"
    16 = "35 Minutos"
    18 = "1 Hora e 45 Minutos"
    19 = "1 Hora e 50 Minutos"
    20 = "1 Hora e 40 Minutos"
    22 = "50 Minutos"
    23 = "  1 Hora e 50Minutos"
    24 = " 1 Hora e 40 Minutos"
    25 = " 1 Hora e 30 Minutos"
    27 = "1 Hora 20 Minutos"
    28 = "1 Hora"
    29 = "1 Hora e 10 Minutos"
    30 = "1 Horas e 30 Minutos"
    31 = "2 Horas e 15 Minutos"
    33 = "1 Hora e 10 Minutos"
    34 = "2 Horas"
    35 = "1 Hora e 40 Minutos"
    36 = "1 Hora e 30 Minutos"
    37 = "1 Hora"
    38 = "1 Hora e 40 Minutos"
    39 = "1 Hora e 30 Minutos"
}

foreach ($row in $values.Keys) {
    $ws.Cells.Item($row, 6).Value = $values[$row]
}

# ---------------------------------------------------------------------
# 6) Column widths: new F (RESULTADO EM HORAS) is wider; G keeps the
#    width the STATUS column already had before the insert.
# ---------------------------------------------------------------------
$ws.Columns.Item(6).ColumnWidth = 25.1
